# Update the "想去人数" (F column) counts across all sheets to reflect
# the newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 73
$ws.Range("F3").Value = 195
$ws.Range("F5").Value = 1623
$ws.Range("F6").Value = 3230
$ws.Range("F7").Value = 768
$ws.Range("F8").Value = 1984
$ws.Range("F9").Value = 1901
$ws.Range("F10").Value = 975
$ws.Range("F11").Value = 340
$ws.Range("F12").Value = 14
$ws.Range("F13").Value = 1574
$ws.Range("F14").Value = 335
$ws.Range("F16").Value = 58
$ws.Range("F17").Value = 1390
$ws.Range("F18").Value = 485
$ws.Range("F19").Value = 596
$ws.Range("F20").Value = 285
$ws.Range("F21").Value = 10454
$ws.Range("F22").Value = 9637
$ws.Range("F23").Value = 833
$ws.Range("F25").Value = 1805
$ws.Range("F27").Value = 371

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 6
$ws.Range("F6").Value = 30

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 56

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 73
$ws.Range("F3").Value = 56
$ws.Range("F4").Value = 195
$ws.Range("F7").Value = 1623
$ws.Range("F8").Value = 3230
$ws.Range("F9").Value = 768
$ws.Range("F10").Value = 1984
$ws.Range("F11").Value = 1901
$ws.Range("F12").Value = 975
$ws.Range("F13").Value = 340
$ws.Range("F14").Value = 14
$ws.Range("F15").Value = 1574
$ws.Range("F16").Value = 335
$ws.Range("F19").Value = 58
$ws.Range("F21").Value = 1390
$ws.Range("F22").Value = 485
$ws.Range("F23").Value = 596
$ws.Range("F24").Value = 285
$ws.Range("F25").Value = 10455
$ws.Range("F26").Value = 9637
$ws.Range("F27").Value = 833
$ws.Range("F29").Value = 1805
$ws.Range("F31").Value = 30
$ws.Range("F33").Value = 371
